$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1262
$ws.Range("J17").Value = 1700
$ws.Range("L17").Value = 5100
$ws.Range("N17").Value = -5436

# Row 64
$ws.Range("H64").Value = 12845.77
$ws.Range("I64").Value = 9999.5
$ws.Range("J64").Value = 17399.8
$ws.Range("K64").Value = 9999.5
$ws.Range("L64").Value = 17399.8
$ws.Range("M64").Value = -9751.5
$ws.Range("N64").Value = -17895.8

# Row 67
$ws.Range("H67").Value = 12845.77
$ws.Range("I67").Value = 9999.5
$ws.Range("J67").Value = 17399.8
$ws.Range("K67").Value = 9999.5
$ws.Range("L67").Value = 17399.8
$ws.Range("M67").Value = -9141.5
$ws.Range("N67").Value = -19115.8

# Row 105
$ws.Range("H105").Value = 17499.5
$ws.Range("J105").Value = 17499.5
$ws.Range("L105").Value = 17499.5
$ws.Range("N105").Value = -24487.5

# Row 125
$ws.Range("H125").Value = 62504944
$ws.Range("I125").Value = 125000376
$ws.Range("K125").Value = 1125003384
$ws.Range("M125").Value = -1125000924

# Row 137
$ws.Range("H137").Value = 1950.1578
$ws.Range("I137").Value = 822.8182
$ws.Range("J137").Value = 3500.25
$ws.Range("K137").Value = 2468.4546
$ws.Range("L137").Value = 10500.75
$ws.Range("M137").Value = 81.54539999999997
$ws.Range("N137").Value = -15600.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 399
$ws.Range("I2").Value = 401.1111
$ws.Range("K2").Value = 401.1111
$ws.Range("M2").Value = -288.1111

# Row 32
$ws.Range("H32").Value = 4295.451
$ws.Range("I32").Value = 3246.449
$ws.Range("J32").Value = 29996
$ws.Range("K32").Value = 3246.449
$ws.Range("L32").Value = 29996
$ws.Range("M32").Value = -2959.449
$ws.Range("N32").Value = -30570

# Row 45
$ws.Range("H45").Value = 1985.1428
$ws.Range("I45").Value = 1979.8
$ws.Range("K45").Value = 1979.8
$ws.Range("M45").Value = -1602.8

# Row 61
$ws.Range("H61").Value = 2218.3333
$ws.Range("I61").Value = 2218.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2218.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2006.3333
$ws.Range("N61").Value = ""

# Row 74
$ws.Range("H74").Value = 1112.25
$ws.Range("I74").Value = 1056.8572
$ws.Range("K74").Value = 1056.8572
$ws.Range("M74").Value = -182.8571999999999

# Row 77
$ws.Range("H77").Value = 1112.25
$ws.Range("I77").Value = 1056.8572
$ws.Range("K77").Value = 5284.286
$ws.Range("M77").Value = -916.2860000000001

# Row 98
$ws.Range("H98").Value = 20797.4
$ws.Range("J98").Value = 20797.4
$ws.Range("L98").Value = 20797.4
$ws.Range("N98").Value = -26787.4

# Row 106
$ws.Range("H106").Value = 22886.334
$ws.Range("J106").Value = 22886.334
$ws.Range("L106").Value = 22886.334
$ws.Range("N106").Value = -25410.334

# Row 116
$ws.Range("H116").Value = 399
$ws.Range("I116").Value = 401.1111
$ws.Range("K116").Value = 401.1111
$ws.Range("M116").Value = 1892.8889

# Row 122
$ws.Range("H122").Value = 1695.3
$ws.Range("I122").Value = 1695.3
$ws.Range("K122").Value = 5085.9
$ws.Range("M122").Value = -2635.9

# Row 132
$ws.Range("H132").Value = 2959.5
$ws.Range("I132").Value = 1866.6666
$ws.Range("K132").Value = 5599.9998
$ws.Range("M132").Value = -3069.9998

# Row 136
$ws.Range("H136").Value = 2218.3333
$ws.Range("I136").Value = 2218.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6654.999899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4104.999899999999
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 399
$ws.Range("I3").Value = 401.1111
$ws.Range("K3").Value = 401.1111
$ws.Range("M3").Value = -287.1111

# Row 100
$ws.Range("H100").Value = 8699
$ws.Range("J100").Value = 8699
$ws.Range("L100").Value = 8699
$ws.Range("N100").Value = -10863

# Row 134
$ws.Range("H134").Value = 2790.2222
$ws.Range("I134").Value = 2538.4285
$ws.Range("K134").Value = 7615.2855
$ws.Range("M134").Value = -5080.2855

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 4157
$ws.Range("J22").Value = 6315
$ws.Range("L22").Value = 6315
$ws.Range("N22").Value = -7015

# Row 106
$ws.Range("H106").Value = 39666.332
$ws.Range("J106").Value = 39666.332
$ws.Range("L106").Value = 39666.332
$ws.Range("N106").Value = -42190.332

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 600.4286
$ws.Range("I5").Value = 301.75
$ws.Range("K5").Value = 905.25
$ws.Range("M5").Value = -793.25

# Row 38
$ws.Range("H38").Value = 924.3333
$ws.Range("J38").Value = 175
$ws.Range("L38").Value = 525
$ws.Range("N38").Value = -1219

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = ""

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = ""

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = ""

# Row 135
$ws.Range("H135").Value = 600.4286
$ws.Range("I135").Value = 301.75
$ws.Range("K135").Value = 2715.75
$ws.Range("M135").Value = -180.75

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""

# Row 44
$ws.Range("H44").Value = 15000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 15000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -16192
$ws.Range("M44").Value = ""

# Row 134
$ws.Range("H134").Value = 59996.332
$ws.Range("J134").Value = 59996.332
$ws.Range("L134").Value = 179988.996
$ws.Range("N134").Value = -185058.996

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2686.8845
$ws.Range("I46").Value = 1910.5555
$ws.Range("J46").Value = 3097.8823
$ws.Range("K46").Value = 1910.5555
$ws.Range("L46").Value = 3097.8823
$ws.Range("M46").Value = -1722.5555
$ws.Range("N46").Value = -3473.8823

# Row 61
$ws.Range("H61").Value = 1864.2778
$ws.Range("I61").Value = 1295.1666
$ws.Range("K61").Value = 1295.1666
$ws.Range("M61").Value = -1093.1666

# Row 63
$ws.Range("H63").Value = 41945
$ws.Range("I63").Value = 41945
$ws.Range("K63").Value = 41945
$ws.Range("M63").Value = -41196

# Row 66
$ws.Range("H66").Value = 41945
$ws.Range("I66").Value = 41945
$ws.Range("K66").Value = 125835
$ws.Range("M66").Value = -122091

# Row 92
$ws.Range("H92").Value = 24000
$ws.Range("J92").Value = 24000
$ws.Range("L92").Value = 24000
$ws.Range("N92").Value = -28992

# Row 113
$ws.Range("H113").Value = 1864.2778
$ws.Range("I113").Value = 1295.1666
$ws.Range("K113").Value = 1295.1666
$ws.Range("M113").Value = 874.8334

# Row 136
$ws.Range("H136").Value = 4987.25
$ws.Range("I136").Value = 4666.3335
$ws.Range("K136").Value = 13999.0005
$ws.Range("M136").Value = -11449.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 29999
$ws.Range("J64").Value = 29999
$ws.Range("L64").Value = 29999
$ws.Range("N64").Value = -30495

# Row 67
$ws.Range("H67").Value = 29999
$ws.Range("J67").Value = 29999
$ws.Range("L67").Value = 29999
$ws.Range("N67").Value = -31715

# Row 132
$ws.Range("H132").Value = 2749.8333
$ws.Range("I132").Value = 3058.8
$ws.Range("J132").Value = 1205
$ws.Range("K132").Value = 9176.400000000001
$ws.Range("L132").Value = 3615
$ws.Range("M132").Value = -6646.400000000001
$ws.Range("N132").Value = -8675

